$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.436.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.80%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.275.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.70%  "

# Row 4
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "123.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.63%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "266.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.16%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.640"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.21%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.624"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.51%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.53%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0944"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.24%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.27%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.18%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.72%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.902"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.41%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.619.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.63%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.277.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.45%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.551.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.25%  "

# Row 19
$ws.Range("E19").Value = "  +0.19%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.82%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.42%  "

# Row 22
$ws.Range("E22").Value = "  -0.76%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.80%  "

# Row 24
$ws.Range("E24").Value = "  -2.13%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.80%  "

# Row 26
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.66%  "

# Row 27
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.76%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.05%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.45%  "

# Row 30
$ws.Range("E30").Value = "  +0.76%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.90%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.30%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0915"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.42%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.09%  "

# Row 35
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.130"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.43%  "

# Row 36
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.49%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0376"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.41%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.31%  "

# Row 39
$ws.Range("E39").Value = "  -2.44%  "

# Row 40
$ws.Range("E40").Value = "  +4.41%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.47%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.01%  "

# Row 43
$ws.Range("E43").Value = "  -2.12%  "

# Row 44
$ws.Range("E44").Value = "  -0.11%  "

# Row 45
$ws.Range("E45").Value = "  -1.55%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -11.77%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "74.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +37.05%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.47%  "

# Row 49
$ws.Range("E49").Value = "  -0.29%  "

# Row 50
$ws.Range("E50").Value = "  +0.17%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "101.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.43%  "
